$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff: A1 1 -> 100, A2 2 -> 102
$ws.Range("A1").Value = 100
$ws.Range("A2").Value = 102

# Reset the active cell selection back to A1 (removes the stray <selection activeCell="E16" sqref="E16"/>)
$ws.Range("A1").Select()
